$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VENTA MENSUAL")

# Mapping of row -> new value for column G (PRESUPUESTO)
$updates = @{
    2  = 0
    5  = 8200
    6  = 3000
    11 = 3500
    12 = 2500
    19 = 3500
    21 = 3000
    22 = 3500
    24 = 3500
    26 = 2000
    27 = 8000
    28 = 3500
    29 = 3500
    30 = 1000
    31 = 3500
    34 = 3500
    35 = 3500
    36 = 5501
    37 = 5800
    39 = 4000
    41 = 3000
    43 = 5800
    44 = 4000
    45 = 3000
    47 = 2300
    48 = 2500
    49 = 3200
    50 = 4300
    52 = 500
    53 = 2100
    55 = 108201
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 7).Value = $updates[$row]
}
